$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 992.1905
$ws.Range("I17").Value = 595.1111
$ws.Range("J17").Value = 1290
$ws.Range("K17").Value = 1785.3333
$ws.Range("L17").Value = 3870
$ws.Range("M17").Value = -1617.3333
$ws.Range("N17").Value = -4206
$ws.Range("H19").Value = 11894.777
$ws.Range("I19").Value = 688.8889
$ws.Range("J19").Value = 23100.666
$ws.Range("K19").Value = 688.8889
$ws.Range("L19").Value = 23100.666
$ws.Range("M19").Value = -513.8889
$ws.Range("N19").Value = -23450.666
$ws.Range("H44").Value = 10700
$ws.Range("J44").Value = 10700
$ws.Range("L44").Value = 10700
$ws.Range("N44").Value = -11624
$ws.Range("H62").Value = 6646.4116
$ws.Range("I62").Value = 3170.7144
$ws.Range("J62").Value = 22866.334
$ws.Range("K62").Value = 3170.7144
$ws.Range("L62").Value = 22866.334
$ws.Range("M62").Value = -2546.7144
$ws.Range("N62").Value = -24114.334
$ws.Range("H65").Value = 6646.4116
$ws.Range("I65").Value = 3170.7144
$ws.Range("J65").Value = 22866.334
$ws.Range("K65").Value = 15853.572
$ws.Range("L65").Value = 114331.67
$ws.Range("M65").Value = -12733.572
$ws.Range("N65").Value = -120571.67
$ws.Range("H86").Value = 69582110
$ws.Range("I86").Value = 84230180
$ws.Range("J86").Value = 3800
$ws.Range("K86").Value = 84230180
$ws.Range("L86").Value = 3800
$ws.Range("M86").Value = -84229057
$ws.Range("N86").Value = -6046
$ws.Range("H89").Value = 69582110
$ws.Range("I89").Value = 84230180
$ws.Range("J89").Value = 3800
$ws.Range("K89").Value = 421150900
$ws.Range("L89").Value = 19000
$ws.Range("M89").Value = -421145284
$ws.Range("N89").Value = -30232
$ws.Range("H112").Value = 10622.223
$ws.Range("J112").Value = 10622.223
$ws.Range("L112").Value = 31866.669
$ws.Range("N112").Value = -34082.669
$ws.Range("H125").Value = 1817
$ws.Range("J125").Value = 1817
$ws.Range("L125").Value = 16353
$ws.Range("N125").Value = -21273
$ws.Range("H128").Value = 46000
$ws.Range("J128").Value = 46000
$ws.Range("L128").Value = 46000
$ws.Range("N128").Value = -55960
$ws.Range("H129").Value = 1582.875
$ws.Range("J129").Value = 2182.238
$ws.Range("L129").Value = 6546.714
$ws.Range("N129").Value = -16546.714
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6339647.5
$ws.Range("I32").Value = 7584660.5
$ws.Range("K32").Value = 7584660.5
$ws.Range("M32").Value = -7584373.5
$ws.Range("H92").Value = 65600
$ws.Range("J92").Value = 65600
$ws.Range("L92").Value = 65600
$ws.Range("N92").Value = -70592
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H134").Value = 2744.0222
$ws.Range("I134").Value = 2696.4243
$ws.Range("K134").Value = 8089.2729
$ws.Range("M134").Value = -5554.2729
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5137.3535
$ws.Range("I31").Value = 1643.9688
$ws.Range("J31").Value = 7373.12
$ws.Range("K31").Value = 1643.9688
$ws.Range("L31").Value = 7373.12
$ws.Range("M31").Value = -1348.9688
$ws.Range("N31").Value = -7963.12
$ws.Range("H34").Value = 5137.3535
$ws.Range("I34").Value = 1643.9688
$ws.Range("J34").Value = 7373.12
$ws.Range("K34").Value = 1643.9688
$ws.Range("L34").Value = 7373.12
$ws.Range("M34").Value = -1441.9688
$ws.Range("N34").Value = -7777.12
$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 641.94934
$ws.Range("I5").Value = 502.09805
$ws.Range("J5").Value = 896.6786
$ws.Range("K5").Value = 1506.29415
$ws.Range("L5").Value = 2690.0358
$ws.Range("M5").Value = -1394.29415
$ws.Range("N5").Value = -2914.0358
$ws.Range("H35").Value = 2256.2222
$ws.Range("I35").Value = 300
$ws.Range("J35").Value = 2815.1428
$ws.Range("K35").Value = 900
$ws.Range("L35").Value = 8445.428400000001
$ws.Range("M35").Value = -612
$ws.Range("N35").Value = -9021.428400000001
$ws.Range("H93").Value = 4965.8335
$ws.Range("I93").Value = 4921
$ws.Range("J93").Value = 4968.4707
$ws.Range("K93").Value = 14763
$ws.Range("L93").Value = 14905.4121
$ws.Range("M93").Value = -12891
$ws.Range("N93").Value = -18649.4121
$ws.Range("H97").Value = 1900
$ws.Range("J97").Value = 1900
$ws.Range("L97").Value = 5700
$ws.Range("N97").Value = -6692
$ws.Range("H113").Value = 886.6667
$ws.Range("I113").Value = 886.6667
$ws.Range("K113").Value = 2660.0001
$ws.Range("M113").Value = -490.0001000000002
$ws.Range("H135").Value = 641.94934
$ws.Range("I135").Value = 502.09805
$ws.Range("J135").Value = 896.6786
$ws.Range("K135").Value = 4518.88245
$ws.Range("L135").Value = 8070.1074
$ws.Range("M135").Value = -1983.88245
$ws.Range("N135").Value = -13140.1074
$ws.Range("H137").Value = 31436.553
$ws.Range("I137").Value = 5944
$ws.Range("J137").Value = 75138.07000000001
$ws.Range("K137").Value = 17832
$ws.Range("L137").Value = 225414.21
$ws.Range("M137").Value = -12732
$ws.Range("N137").Value = -235614.21
$ws.Range("H141").Value = 9395.736999999999
$ws.Range("I141").Value = 8171.9
$ws.Range("K141").Value = 24515.7
$ws.Range("M141").Value = -19335.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 18800
$ws.Range("J58").Value = 16000
$ws.Range("L58").Value = 16000
$ws.Range("N58").Value = -16554
$ws.Range("H80").Value = 29345946
$ws.Range("I80").Value = 35135136
$ws.Range("K80").Value = 35135136
$ws.Range("M80").Value = -35134138
$ws.Range("H83").Value = 29345946
$ws.Range("I83").Value = 35135136
$ws.Range("K83").Value = 175675680
$ws.Range("M83").Value = -175670688
$ws.Range("H102").Value = 1422.9286
$ws.Range("I102").Value = 1460.0834
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1460.0834
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 161.9166
$ws.Range("N102").Value = -4444
$ws.Range("H116").Value = 38566.168
$ws.Range("J116").Value = 38566.168
$ws.Range("L116").Value = 38566.168
$ws.Range("N116").Value = -47744.168
$ws.Range("H126").Value = 3427.6667
$ws.Range("I126").Value = 3350.5
$ws.Range("J126").Value = 3455.7273
$ws.Range("K126").Value = 10051.5
$ws.Range("L126").Value = 10367.1819
$ws.Range("M126").Value = -7581.5
$ws.Range("N126").Value = -15307.1819
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2984.0981
$ws.Range("I136").Value = 2609.3823
$ws.Range("J136").Value = 3733.5293
$ws.Range("K136").Value = 7828.146900000001
$ws.Range("L136").Value = 11200.5879
$ws.Range("M136").Value = -5278.146900000001
$ws.Range("N136").Value = -16300.5879
